function ConvertTo-BGRLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Re-colour the deck's theme from "Integral" to the stock "Office Theme"
# palette (Design > Themes > Office Theme), applied through the first
# slide's ThemeColorScheme, which is shared by every slide that hangs off
# the single Slide Master in this deck.
$officeThemeColors = @(
    "000000",  # Dark 1
    "FFFFFF",  # Light 1
    "44546A",  # Dark 2
    "E7E6E6",  # Light 2
    "5B9BD5",  # Accent 1
    "ED7D31",  # Accent 2
    "A5A5A5",  # Accent 3
    "FFC000",  # Accent 4
    "4472C4",  # Accent 5
    "70AD47",  # Accent 6
    "0563C1",  # Hyperlink
    "954F72"   # Followed Hyperlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-BGRLong $officeThemeColors[$i - 1]
}
